# Fill in the AvailabilityZone (column D) values for the EC2/RDS rows so the
# Jinja2 templates have a region to render for each record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "us-east-1"
$ws.Range("D3").Value = "us-west-1"
$ws.Range("D4").Value = "us-central-2"
